$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 04:54"

# Row 26 - Pakistan
$ws.Range("B26").Value = 324744
$ws.Range("C26").Value = 667
$ws.Range("D26").Value = 308674
$ws.Range("E26").Value = 9378
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 6692

# Row 30 - Belgica
$ws.Range("B30").Value = 240159
$ws.Range("C30").Value = 9679
$ws.Range("D30").Value = 21476
$ws.Range("E30").Value = 208194
$ws.Range("G30").Value = 46
$ws.Range("H30").Value = 10489

# Row 45 - Kazajistan
$ws.Range("B45").Value = 109766
$ws.Range("C45").Value = 143
$ws.Range("D45").Value = 105301
$ws.Range("E45").Value = 2669
$ws.Range("H45").Value = 1796

# Row 53 - Honduras
$ws.Range("B53").Value = 90232
$ws.Range("C53").Value = 851
$ws.Range("D53").Value = 35930
$ws.Range("E53").Value = 51720
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 2582
